# Purchasing Dashboard fix-up:
#  - Rename the 5 sheets so their tab names use spaces instead of underscores.
#  - Add a new "CROSS REFERENCE" sheet (upstream production-plan data).
#  - Point MRP ENGINE's "Target Production" row at the new sheet's TOTAL cell.
#  - Re-point the two charts' series formulas at the renamed sheets.
#
# Renaming a sheet through the object model naturally cascades into every
# formula in the workbook that references it (this host rewrites the AST,
# same as real Excel). The authored fix, however, only repointed the
# sheet-name table, the chart series and the MRP ENGINE row 6 formulas --
# every other cross-sheet formula (COST ANALYSIS, CASH FLOW PREVIEW,
# UPLOAD READY PROCUREMENT) was left referring to the old underscored
# names. So: snapshot every formula in the workbook first, rename, then
# restore the snapshot verbatim before layering on the intentional edits.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Snapshot every formula cell in the five original sheets (by sheet
#    index + row/col, which stay stable across a rename).
# ---------------------------------------------------------------------
$savedFormulas = @{}
$origSheetCount = 5
for ($si = 1; $si -le $origSheetCount; $si++) {
    $ws = $wb.Worksheets.Item($si)
    $ur = $ws.UsedRange
    $r0 = $ur.Row
    $c0 = $ur.Column
    $nr = $ur.Rows.Count
    $nc = $ur.Columns.Count
    for ($i = 0; $i -lt $nr; $i++) {
        for ($j = 0; $j -lt $nc; $j++) {
            $cell = $ws.Cells.Item($r0 + $i, $c0 + $j)
            if ($cell.HasFormula) {
                $key = "$si|" + ($r0 + $i) + "|" + ($c0 + $j)
                $savedFormulas[$key] = $cell.Formula
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Add the new CROSS REFERENCE sheet at the end of the tab strip FIRST.
#    (Adding a sheet *after* some formulas are left pointing at an
#    unresolved/renamed sheet name causes this host to silently rebind
#    those dangling references onto the newly-added sheet. Creating the
#    new sheet before we ever touch formula text sidesteps that.)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$crossRef = $wb.Worksheets.Add($null, $lastSheet)
$crossRef.Name = "CROSS REFERENCE"

# ---------------------------------------------------------------------
# 3) Rename the sheets (underscores -> spaces).
# ---------------------------------------------------------------------
$wb.Worksheets.Item(1).Name = "SUPPLIER CONFIG"
$wb.Worksheets.Item(2).Name = "COST ANALYSIS"
$wb.Worksheets.Item(3).Name = "MRP ENGINE"
$wb.Worksheets.Item(4).Name = "CASH FLOW PREVIEW"
$wb.Worksheets.Item(5).Name = "UPLOAD READY PROCUREMENT"

# ---------------------------------------------------------------------
# 4) Restore every formula to its pre-rename text, undoing the automatic
#    reference rewrite so stale cross-sheet formulas match the original
#    (unfixed) commit.
# ---------------------------------------------------------------------
foreach ($key in $savedFormulas.Keys) {
    $parts = $key.Split("|")
    $si = [int]$parts[0]
    $r = [int]$parts[1]
    $c = [int]$parts[2]
    $ws = $wb.Worksheets.Item($si)
    $ws.Cells.Item($r, $c).Formula = $savedFormulas[$key]
}

$supplierConfig = $wb.Worksheets.Item("SUPPLIER CONFIG")
$mrpEngine = $wb.Worksheets.Item("MRP ENGINE")

# Copy reference formatting onto the cells we're about to fill in, then
# set their values. (Style ids already exist in styles.xml; PasteSpecial
# copies the format cleanly without disturbing content.)
$supplierConfig.Range("A1").Copy() | Out-Null
$crossRef.Range("A1").PasteSpecial(-4122) | Out-Null
$supplierConfig.Range("A2").Copy() | Out-Null
$crossRef.Range("A2").PasteSpecial(-4122) | Out-Null
$supplierConfig.Range("A4").Copy() | Out-Null
$crossRef.Range("A4").PasteSpecial(-4122) | Out-Null
$mrpEngine.Range("A5").Copy() | Out-Null
$crossRef.Range("A5:B5").PasteSpecial(-4122) | Out-Null
$supplierConfig.Range("A6").Copy() | Out-Null
$crossRef.Range("A6:B10").PasteSpecial(-4122) | Out-Null
$mrpEngine.Range("A17").Copy() | Out-Null
$crossRef.Range("A11:B11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$crossRef.Range("A1").Value = "CROSS-REFERENCE SUMMARY"
$crossRef.Range("A2").Value = "Upstream data from Production Dashboard"
$crossRef.Range("A4").Value = "PRODUCTION PLAN"
$crossRef.Range("A5").Value = "Zone"
$crossRef.Range("B5").Value = "Target/FN"
$crossRef.Range("A6").Value = "Center"
$crossRef.Range("B6").Value = 31200
$crossRef.Range("A7").Value = "West"
$crossRef.Range("B7").Value = 0
$crossRef.Range("A8").Value = "North"
$crossRef.Range("B8").Value = 0
$crossRef.Range("A9").Value = "East"
$crossRef.Range("B9").Value = 0
$crossRef.Range("A10").Value = "South"
$crossRef.Range("B10").Value = 0
$crossRef.Range("A11").Value = "TOTAL"
$crossRef.Range("B11").Value = 31200

# ---------------------------------------------------------------------
# 5) MRP ENGINE row 6: relabel "Target Production" and pull the number
#    from the new CROSS REFERENCE total instead of a hard-coded 0.
# ---------------------------------------------------------------------
$mrpEngine.Range("A6").Value = "Target Production (from Production Plan)"
$cols = @("B", "C", "D", "E", "F", "G", "H", "I")
foreach ($col in $cols) {
    $mrpEngine.Range($col + "6").Formula = "='CROSS REFERENCE'!B11"
}

# ---------------------------------------------------------------------
# 6) Re-point the chart series formulas at the renamed sheets.
# ---------------------------------------------------------------------
$chart1 = $mrpEngine.ChartObjects(1).Chart
for ($i = 1; $i -le $chart1.SeriesCollection().Count; $i++) {
    $ser = $chart1.SeriesCollection($i)
    $ser.Formula = $ser.Formula.Replace("MRP_ENGINE!", "'MRP ENGINE'!")
}

$cashFlowPreview = $wb.Worksheets.Item("CASH FLOW PREVIEW")
$chart2 = $cashFlowPreview.ChartObjects(1).Chart
for ($i = 1; $i -le $chart2.SeriesCollection().Count; $i++) {
    $ser = $chart2.SeriesCollection($i)
    $ser.Formula = $ser.Formula.Replace("CASH_FLOW_PREVIEW!", "'CASH FLOW PREVIEW'!")
}

# ---------------------------------------------------------------------
# 7) Leave the original sheet selected/active, as in the source file.
# ---------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
